$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 125: Preset Games entry
$ws.Range("A125").Value = 41994
$ws.Range("B125").Value = 0.0083333333333333332
$ws.Range("C125").Value = 0.024305555555555556
$ws.Range("D125").Value = 5
$ws.Range("E125").Formula = "=IF(AND(NOT(ISBLANK(B125)),NOT(ISBLANK(C125))), (C125-B125) * 24 - D125/60, """")"
$ws.Range("F125").Value = "Preset Games"

# Row 126: Custom Games entry
$ws.Range("A126").Value = 41996
$ws.Range("B126").Value = 0.56527777777777777
$ws.Range("C126").Value = 0.61736111111111114
$ws.Range("D126").Value = 10
$ws.Range("E126").Formula = "=IF(AND(NOT(ISBLANK(B126)),NOT(ISBLANK(C126))), (C126-B126) * 24 - D126/60, """")"
$ws.Range("F126").Value = "Custom Games"

# Update the selection / active cell to match the author's cursor position
$ws.Range("A127").Select()

$excel.CalculateFullRebuild()
